$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Update Price (column D) values per diff
Set-TextValue $ws.Range("D2") "247.63"
Set-TextValue $ws.Range("D3") "21.64"
Set-TextValue $ws.Range("D4") "5.283"
Set-TextValue $ws.Range("D5") "0.05588"
Set-TextValue $ws.Range("D6") "3.397"
Set-TextValue $ws.Range("D8") "0.8149"
Set-TextValue $ws.Range("D9") "0.9656"
Set-TextValue $ws.Range("D10") "0.1408"
Set-TextValue $ws.Range("D11") "0.07409"
Set-TextValue $ws.Range("D12") "0.03142"
Set-TextValue $ws.Range("D13") "0.03034"
Set-TextValue $ws.Range("D14") "0.09306"
Set-TextValue $ws.Range("D15") "3.550"
Set-TextValue $ws.Range("D16") "0.001620"
Set-TextValue $ws.Range("D17") "0.04700"
Set-TextValue $ws.Range("D18") "0.0005770"
Set-TextValue $ws.Range("D19") "0.006377"
Set-TextValue $ws.Range("D20") "0.005023"
Set-TextValue $ws.Range("D22") "0.0001498"
Set-TextValue $ws.Range("D23") "3.744"
Set-TextValue $ws.Range("D24") "2.123"
Set-TextValue $ws.Range("D25") "0.3252"
Set-TextValue $ws.Range("D26") "0.1249"
Set-TextValue $ws.Range("D28") "0.0003100"
Set-TextValue $ws.Range("D40") "0.03927"
Set-TextValue $ws.Range("D41") "0.007071"
Set-TextValue $ws.Range("D42") "0.1050"
Set-TextValue $ws.Range("D43") "0.003067"
Set-TextValue $ws.Range("D44") "0.007834"
Set-TextValue $ws.Range("D45") "0.00005814"
Set-TextValue $ws.Range("D46") "0.00000000750"
Set-TextValue $ws.Range("D47") "0.0005500"
Set-TextValue $ws.Range("D48") "0.6800"
Set-TextValue $ws.Range("D49") "0.1517"
Set-TextValue $ws.Range("D50") "0.00002100"
Set-TextValue $ws.Range("D51") "0.01010"

# Update Hora (column G) values from 3 to 4 for all data rows 2-51
for ($r = 2; $r -le 51; $r++) {
    Set-TextValue $ws.Cells.Item($r, 7) "4"
}
